$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update frequency counts for the "same-category" rows (3-6)
$ws.Range("C3").Value = 466
$ws.Range("C4").Value = 396
$ws.Range("C5").Value = 297
$ws.Range("C6").Value = 228

# Row 7: Kids & Toys -> Storage & Organization, frequency 132
$ws.Range("A7").Value = "Kids & Toys"
$ws.Range("B7").Value = "Storage & Organization"
$ws.Range("C7").Value = 132

# Row 8: Kids & Toys -> Kids & Toys, frequency 110
$ws.Range("A8").Value = "Kids & Toys"
$ws.Range("B8").Value = "Kids & Toys"
$ws.Range("C8").Value = 110

# Row 9: Fashion & Accessories -> Storage & Organization, frequency 89
$ws.Range("A9").Value = "Fashion & Accessories"
$ws.Range("B9").Value = "Storage & Organization"
$ws.Range("C9").Value = 89

# Row 10: Storage & Organization -> Fashion & Accessories, frequency 75
$ws.Range("A10").Value = "Storage & Organization"
$ws.Range("B10").Value = "Fashion & Accessories"
$ws.Range("C10").Value = 75

# Row 11: Storage & Organization -> Kids & Toys, frequency 70
$ws.Range("A11").Value = "Storage & Organization"
$ws.Range("B11").Value = "Kids & Toys"
$ws.Range("C11").Value = 70
